$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (C++ / A7): Method Overloading completed (2) and mark "completed" green box
$ws.Range("G7").Copy()
$ws.Range("H7").PasteSpecial(-4122)
$ws.Range("H7").Value = 2

$ws.Range("I6").Copy()
$ws.Range("I7").PasteSpecial(-4122)

# Row 8 (Python / A8): update Smalltalk/Garbage Collection value, add Method Overloading, mark completed
$ws.Range("B8").Value = 2

$ws.Range("G8").Copy()
$ws.Range("H8").PasteSpecial(-4122)
$ws.Range("H8").Value = 0

$ws.Range("I6").Copy()
$ws.Range("I8").PasteSpecial(-4122)

# Row 9 (Perl / A9): update values, add Method Overloading, mark completed
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 2
$ws.Range("E9").Value = 2

$ws.Range("G9").Copy()
$ws.Range("H9").PasteSpecial(-4122)
$ws.Range("H9").Value = 0

$ws.Range("I6").Copy()
$ws.Range("I9").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Move the active selection to I10
$ws.Range("I10").Select()
